$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 21:34"

# --- Updated case counters for existing countries (values refreshed) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1338718
$ws.Range("C4").Value = 16933
$ws.Range("D4").Value = 225004
$ws.Range("E4").Value = 1034115
$ws.Range("G4").Value = 984
$ws.Range("H4").Value = 79599

# Row 9: Francia
$ws.Range("B9").Value = 176658
$ws.Range("C9").Value = 579
$ws.Range("E9").Value = 94310

# Row 10: Alemania
$ws.Range("B10").Value = 171145
$ws.Range("C10").Value = 557
$ws.Range("E10").Value = 20313
$ws.Range("G10").Value = 22
$ws.Range("H10").Value = 7532

# Row 106: El Salvador
$ws.Range("E106").Value = 491
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 17

# Row 107: Costa Rica
$ws.Range("B107").Value = 780
$ws.Range("C107").Value = 7
$ws.Range("D107").Value = 480
$ws.Range("E107").Value = 294

# --- Guinea-Bisau moves into the sorted table at row 116, pushing San
#     Marino / Georgia / Gabon / Tayikistan down one row each (their data
#     shifts with them; row 121 Jordania is unaffected) ---

# Row 116: now Guinea-Bisau (new data)
$ws.Range("A116").Value = "Guinea-Bisau"
$ws.Range("B116").Value = 641
$ws.Range("C116").Value = 47
$ws.Range("D116").Value = 25
$ws.Range("E116").Value = 613
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 3

# Row 117: now San Marino (previous row-116 data)
$ws.Range("A117").Value = "San Marino"
$ws.Range("B117").Value = 637
$ws.Range("C117").Value = 14
$ws.Range("D117").Value = 126
$ws.Range("E117").Value = 470
$ws.Range("F117").Value = 3
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 41

# Row 118: now Georgia (previous row-117 data)
$ws.Range("A118").Value = "Georgia"
$ws.Range("B118").Value = 626
$ws.Range("C118").Value = 3
$ws.Range("D118").Value = 297
$ws.Range("E118").Value = 319
$ws.Range("F118").Value = 6
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 10

# Row 119: now Gabon (previous row-118 data)
$ws.Range("A119").Value = "Gabon"
$ws.Range("B119").Value = 620
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 110
$ws.Range("E119").Value = 502
$ws.Range("F119").Value = 1
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 8

# Row 120: now Tayikistan (previous row-119 data)
$ws.Range("A120").Value = "Tayikistan"
$ws.Range("B120").Value = 612
$ws.Range("C120").Value = 90
$ws.Range("D120").Value = 0
$ws.Range("E120").Value = 592
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 8
$ws.Range("H120").Value = 20

# --- Butan / Islas Virgenes Britanicas swap places at rows 212-213 ---

# Row 212: now Butan (previous row-213 D/H values; rest identical)
$ws.Range("A212").Value = "Butan"
$ws.Range("D212").Value = 5
$ws.Range("H212").Value = 0

# Row 213: now Islas Virgenes Britanicas (previous row-212 D/H values)
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 4
$ws.Range("H213").Value = 1
